$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-index the food samples so the table offers the first 3 (0,1,2) plus the rest,
# shuffling both the displayed food name and its associated nutritional data per row.

# Row 2 -> Flan (index 3)
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Flan"
$ws.Range("C2").Value = 244
$ws.Range("D2").Value = 7.4
$ws.Range("E2").Value = 3.1
$ws.Range("F2").Value = 37.4
$ws.Range("G2").Value = 33.3
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 0

# Row 3 -> Manzana Asada (index 4)
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Manzana Asada"
$ws.Range("C3").Value = 161
$ws.Range("D3").Value = 5.6
$ws.Range("E3").Value = 3.2
$ws.Range("F3").Value = 25.3
$ws.Range("G3").Value = 19.6
$ws.Range("H3").Value = 0.57
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 0

# Row 4 -> Tarta de queso (index 0)
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "Tarta de queso"
$ws.Range("C4").Value = 547
$ws.Range("D4").Value = 18.9
$ws.Range("E4").Value = 9.4
$ws.Range("F4").Value = 75.9
$ws.Range("G4").Value = 41.6
$ws.Range("H4").Value = 17.5
$ws.Range("I4").Value = 15
$ws.Range("J4").Value = 1

# Row 5 -> Torrijas (index 2)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Torrijas"
$ws.Range("C5").Value = 566
$ws.Range("D5").Value = 39.6
$ws.Range("E5").Value = 8.4
$ws.Range("F5").Value = 38.9
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 12.4
$ws.Range("I5").Value = 15
$ws.Range("J5").Value = 1

# Row 6 -> Tarta de chocolate (index 1)
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Tarta de chocolate"
$ws.Range("C6").Value = 922
$ws.Range("D6").Value = 51.3
$ws.Range("E6").Value = 29.3
$ws.Range("F6").Value = 94.5
$ws.Range("G6").Value = 73.4
$ws.Range("H6").Value = 15.8
$ws.Range("I6").Value = 15
$ws.Range("J6").Value = 2
